$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Website Name", "Profile URL")
    ,@("35photo", "https://35photo.pro/@gorbash/")
    ,@("Audiojungle", "https://audiojungle.net/user/gorbash")
    ,@("Bandlab", "https://www.bandlab.com/api/v1.3/users/gorbash")
    ,@("ask.fm", "https://ask.fm/gorbash")
    ,@("BodyBuilding.com", "http://api.bodybuilding.com/api-proxy/bbc/get?slug=gorbash")
    ,@("cfx.re", "https://forum.cfx.re/u/gorbash.json")
    ,@("chatango.com", "https://gorbash.chatango.com")
    ,@("chaturbate", "https://chaturbate.com/gorbash/")
    ,@("Chess.com", "https://api.chess.com/pub/player/gorbash")
    ,@("cda.pl", "https://www.cda.pl/gorbash")
    ,@("Code Project", "https://www.codeproject.com/Members/gorbash")
    ,@("Chomikuj.pl", "https://chomikuj.pl/gorbash/")
    ,@("Codewars", "https://www.codewars.com/users/gorbash")
    ,@("cHEEZburger", "https://profile.cheezburger.com/gorbash")
    ,@("Chyoa", "https://chyoa.com/user/gorbash")
    ,@("crowdin", "https://crowdin.com/profile/gorbash")
    ,@("Mastodon-Defcon", "https://defcon.social/@gorbash")
    ,@("diigo", "https://www.diigo.com/interact_api/load_profile_info?name=gorbash")
    ,@("Disqus", "https://disqus.com/by/gorbash/")
    ,@("DeviantArt", "https://www.deviantart.com/gorbash")
    ,@("Dating.ru", "https://dating.ru/gorbash/")
    ,@("Discogs", "https://api.discogs.com/users/gorbash")
    ,@("DockerHub", "https://hub.docker.com/v2/users/gorbash/")
    ,@("Cults3D", "https://cults3d.com/en/users/gorbash/creations")
    ,@("datezone", "https://www.datezone.com/users/gorbash/")
    ,@("Duolingo", "https://www.duolingo.com/2017-06-30/users?username=gorbash&_=1628308619574")
    ,@("Etsy", "https://www.etsy.com/people/gorbash")
    ,@("Demotywatory", "https://demotywatory.pl/user/gorbash")
    ,@("Fark", "https://www.fark.com/users/gorbash")
    ,@("Fabswingers", "https://www.fabswingers.com/profile/gorbash")
    ,@("eBay", "https://www.ebay.com/usr/gorbash")
    ,@("FatSecret", "https://www.fatsecret.com/member/gorbash")
    ,@("fansly", "https://apiv2.fansly.com/api/v1/account?usernames=gorbash")
    ,@("Foursquare", "https://foursquare.com/gorbash")
    ,@("Filmweb", "https://www.filmweb.pl/user/gorbash")
    ,@("Flipboard", "https://flipboard.com/@gorbash")
    ,@("FurAffinity", "https://www.furaffinity.net/user/gorbash")
    ,@("Geocaching", "https://www.geocaching.com/p/?u=gorbash")
    ,@("GitHub", "https://github.com/gorbash")
    ,@("Gravatar", "https://en.gravatar.com/gorbash.json")
    ,@("GitLab", "https://gitlab.com/gorbash")
    ,@("giters", "https://giters.com/gorbash")
    ,@("Hacker News", "https://news.ycombinator.com/user?id=gorbash")
    ,@("HudsonRock", "https://cavalier.hudsonrock.com/api/json/v2/osint-tools/search-by-username?username=gorbash")
    ,@("HackerOne", "https://hackerone.com/gorbash")
    ,@("FriendFinder-X", "https://www.friendfinder-x.com/profile/gorbash")
    ,@("IFTTT", "https://ifttt.com/p/gorbash")
    ,@("hugging_face", "https://huggingface.co/gorbash")
    ,@("imgur", "https://api.imgur.com/account/v1/accounts/gorbash?client_id=546c25a59c58ad7&include=trophies%2Cmedallions")
    ,@("inaturalist", "https://inaturalist.nz/people/gorbash")
    ,@("InkBunny", "https://inkbunny.net/gorbash")
    ,@("Instagram2", "https://dumpoir.com/v/gorbash")
    ,@("kaggle", "https://www.kaggle.com/gorbash")
    ,@("issuu", "https://issuu.com/gorbash")
    ,@("Internet Archive User Search", "https://archive.org/advancedsearch.php?q=gorbash&output=json")
    ,@("Keybase", "https://keybase.io/gorbash")
    ,@("Kongregate", "https://www.kongregate.com/accounts/gorbash")
    ,@("Linktree", "https://linktr.ee/gorbash")
    ,@("Jeuxvideo", "https://www.jeuxvideo.com/profil/gorbash?mode=infos")
    ,@("lichess", "https://lichess.org/@/gorbash")
    ,@("LibraryThing", "https://www.librarything.com/profile/gorbash")
    ,@("waytohey", "https://waytohey.com/gorbash")
    ,@("Mastodon API", "https://mastodon.social/api/v2/search?q=gorbash&limit=1&type=accounts")
    ,@("MCUUID (Minecraft)", "https://playerdb.co/api/player/minecraft/gorbash")
    ,@("Livejournal", "https://gorbash.livejournal.com")
    ,@("Mod DB", "https://www.moddb.com/members/gorbash")
    ,@("MySpace", "https://myspace.com/gorbash")
    ,@("MyAnimeList", "https://myanimelist.net/profile/gorbash")
    ,@("OpenStreetMap", "https://www.openstreetmap.org/user/gorbash")
    ,@("MCName (Minecraft)", "https://mcname.info/en/search?q=gorbash")
    ,@("Playstation Network", "https://psnprofiles.com/xhr/search/users?q=gorbash")
    ,@("pikabu", "https://pikabu.ru/@gorbash")
    ,@("Pinterest", "https://www.pinterest.com/gorbash/")
    ,@("Pokerstrategy", "http://www.pokerstrategy.net/user/gorbash/profile/")
    ,@("Periscope", "https://www.periscope.tv/gorbash")
    ,@("palnet", "https://www.palnet.io/@gorbash/")
    ,@("Pornhub Users", "https://www.pornhub.com/users/gorbash")
    ,@("prv.pl", "https://www.prv.pl/osoba/gorbash")
    ,@("public", "https://public.com/@gorbash")
    ,@("Roblox", "https://auth.roblox.com/v1/usernames/validate?username=gorbash&birthday=2019-12-31T23:00:00.000Z")
    ,@("rsi", "https://robertsspaceindustries.com/citizens/gorbash")
    ,@("RumbleUser", "https://rumble.com/user/gorbash")
    ,@("slideshare", "https://www.slideshare.net/gorbash")
    ,@("sentimente", "https://www.sentimente.com/amp/gorbash.html")
    ,@("SFD", "https://www.sfd.pl/profile/gorbash")
    ,@("Snapchat Stories", "https://story.snapchat.com/s/gorbash")
    ,@("SoundCloud", "https://soundcloud.com/gorbash")
    ,@("Sourceforge", "https://sourceforge.net/u/gorbash/profile")
    ,@("sofurry", "https://gorbash.sofurry.com")
    ,@("Spotify", "https://open.spotify.com/user/gorbash")
    ,@("Steam", "https://steamcommunity.com/id/gorbash")
    ,@("tabletoptournament", "https://www.tabletoptournaments.net/eu/player/gorbash")
    ,@("Telegram", "https://t.me/gorbash")
    ,@("themeforest", "https://themeforest.net/user/gorbash")
    ,@("TikTok", "https://www.tiktok.com/oembed?url=https://www.tiktok.com/@gorbash")
    ,@("Trello", "https://trello.com/1/Members/gorbash?fields=activityBlocked%2CavatarUrl%2Cbio%2CbioData%2Cconfirmed%2CfullName%2CidEnterprise%2CidMemberReferrer%2Cinitials%2CmemberType%2CnonPublic%2Cproducts%2Curl%2Cusername")
    ,@("trakt", "https://trakt.tv/users/gorbash")
    ,@("tumblr", "https://gorbash.tumblr.com")
    ,@("Twitch", "https://twitchtracker.com/gorbash")
    ,@("Twitter", "https://nitter.privacydev.net/gorbash")
    ,@("Venmo", "https://account.venmo.com/u/gorbash")
    ,@("Vimeo", "https://vimeo.com/gorbash")
    ,@("untappd", "https://untappd.com/user/gorbash/")
    ,@("wattpad", "https://www.wattpad.com/api/v3/users/gorbash?fields=username%2Cname%2Cdescription%2Cavatar%2CbackgroundUrl%2CcreateDate%2Clocation%2Cfollowing%2CfollowingRequest%2CnumFollowing%2Cfollower%2CfollowerRequest%2CnumFollowers%2CnumLists%2CnumStoriesPublished%2CvotesReceived%2Cfacebook%2Ctwitter%2Cwebsite%2Csmashwords%2Chighlight_colour%2Chtml_enabled%2Cverified%2Cambassador%2Cwattpad_squad%2Cis_staff%2Cprograms(wattpad_stars)%2CisPrivate%2CisMuted%2CexternalId%2Cnotes")
    ,@("Wikidot", "http://www.wikidot.com/user:info/gorbash")
    ,@("Xbox Gamertag", "https://www.xboxgamertag.com/search/gorbash")
    ,@("WordPress Support", "https://wordpress.org/support/users/gorbash/")
    ,@("VK", "https://vk.com/gorbash")
    ,@("xHamster", "https://xhamster.com/users/gorbash")
    ,@("YouTube User", "https://www.youtube.com/user/gorbash/about")
    ,@("YouTube User2", "https://www.youtube.com/@gorbash")
    ,@("XVIDEOS-profiles", "https://www.xvideos.com/profiles/gorbash")
    ,@("XVIDEOS-models", "https://www.xvideos.com/models/gorbash")
    ,@("Zbiornik", "https://mini.zbiornik.com/gorbash")
    ,@("XNXX", "https://www.xnxx.com/mobile/profile/gorbash")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-extra row 117 (previously had data, now should be empty)
$ws.Range("A117:B117").ClearContents()

Write-Host "Done. Rows written: $($data.Count)"